$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove any existing hyperlinks first; they will be re-added for every row below
$ws.Cells.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = "2025-12-08 18:28:28"
$ws.Range("B2").Value = "【完全在宅】ChatGPT・AI活用講師募集|IT/業務支援経験者歓迎!"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5449394"
$ws.Range("G2").Value = 600
$ws.Range("H2").Value = "🔥AI,GPT"
$ws.Hyperlinks.Add($ws.Range("F2"), $ws.Range("F2").Value()) | Out-Null

# Row 3
$ws.Range("A3").Value = "2025-12-08 18:28:28"
$ws.Range("B3").Value = "DreamWeaver – 夢日記 + 睡眠導入 + AI分析のアプリ開発"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5449048"
$ws.Range("G3").Value = 370
$ws.Range("H3").Value = "🔥AI,Ai ◆開発 ◇アプリ"
$ws.Hyperlinks.Add($ws.Range("F3"), $ws.Range("F3").Value()) | Out-Null

# Row 4
$ws.Range("A4").Value = "2025-12-08 18:28:28"
$ws.Range("B4").Value = "B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5445154"
$ws.Range("G4").Value = 368
$ws.Range("H4").Value = "🔥AI,Ai ◆開発"
$ws.Hyperlinks.Add($ws.Range("F4"), $ws.Range("F4").Value()) | Out-Null

# Row 5
$ws.Range("A5").Value = "2025-12-08 18:28:28"
$ws.Range("B5").Value = "法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5445159"
$ws.Range("G5").Value = 368
$ws.Range("H5").Value = "🔥AI,Ai ◆開発"
$ws.Hyperlinks.Add($ws.Range("F5"), $ws.Range("F5").Value()) | Out-Null

# Row 6
$ws.Range("A6").Value = "2025-12-08 18:28:28"
$ws.Range("B6").Value = "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5434128"
$ws.Range("G6").Value = 368
$ws.Range("H6").Value = "🔥AI,Ai ◆開発"
$ws.Hyperlinks.Add($ws.Range("F6"), $ws.Range("F6").Value()) | Out-Null

# Row 7
$ws.Range("A7").Value = "2025-12-08 18:28:28"
$ws.Range("B7").Value = "企業のMicrosoft Copilot導入・活用支援AIコンサルタント募集(研修講師・メンター)"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5434363"
$ws.Range("G7").Value = 348
$ws.Range("H7").Value = "🔥AI,Ai ◆コンサル"
$ws.Hyperlinks.Add($ws.Range("F7"), $ws.Range("F7").Value()) | Out-Null

# Row 8
$ws.Range("A8").Value = "2025-12-08 18:28:28"
$ws.Range("B8").Value = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5427956"
$ws.Range("G8").Value = 310
$ws.Range("H8").Value = "🔥AI,Ai"
$ws.Hyperlinks.Add($ws.Range("F8"), $ws.Range("F8").Value()) | Out-Null

# Row 9
$ws.Range("A9").Value = "2025-12-08 18:28:28"
$ws.Range("B9").Value = "【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5439158"
$ws.Range("G9").Value = 303
$ws.Range("H9").Value = "🔥AI,Ai"
$ws.Hyperlinks.Add($ws.Range("F9"), $ws.Range("F9").Value()) | Out-Null

# Row 10
$ws.Range("A10").Value = "2025-12-08 18:28:28"
$ws.Range("B10").Value = "【TypeScript/Clasp必須】LINE WORKS連携ファイル自動保存システムのGAS開発"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5449466"
$ws.Range("G10").Value = 193
$ws.Range("H10").Value = "🔥TypeScript ◆開発"
$ws.Hyperlinks.Add($ws.Range("F10"), $ws.Range("F10").Value()) | Out-Null

# Row 11
$ws.Range("A11").Value = "2025-12-08 18:28:28"
$ws.Range("B11").Value = "【フリーランス募集】訪問看護向けスケジュール管理アプリ開発"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5449280"
$ws.Range("G11").Value = 135
$ws.Range("H11").Value = "◆開発 ◇アプリ"
$ws.Hyperlinks.Add($ws.Range("F11"), $ws.Range("F11").Value()) | Out-Null

# Row 12
$ws.Range("A12").Value = "2025-12-08 18:28:28"
$ws.Range("B12").Value = "[週2常駐] Laravel + Vue.js 基幹業務システム開発"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5449536"
$ws.Range("G12").Value = 125
$ws.Range("H12").Value = "◆開発,システム開発"
$ws.Hyperlinks.Add($ws.Range("F12"), $ws.Range("F12").Value()) | Out-Null

# Row 13
$ws.Range("A13").Value = "2025-12-08 18:28:28"
$ws.Range("B13").Value = "【急募】紙の伝票をWEBシステムへ自動データ入力開発"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5449142"
$ws.Range("G13").Value = 90
$ws.Range("H13").Value = "◆開発"
$ws.Hyperlinks.Add($ws.Range("F13"), $ws.Range("F13").Value()) | Out-Null

# Row 14
$ws.Range("A14").Value = "2025-12-08 18:28:28"
$ws.Range("B14").Value = "資格試験用の過去問のクイズ型式のiPhone用アプリ開発"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5449723"
$ws.Range("G14").Value = 85
$ws.Range("H14").Value = "◆開発 ◇アプリ"
$ws.Hyperlinks.Add($ws.Range("F14"), $ws.Range("F14").Value()) | Out-Null

# Row 15
$ws.Range("A15").Value = "2025-12-08 18:28:28"
$ws.Range("B15").Value = "自動出品システムの開発"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5449232"
$ws.Range("G15").Value = 83
$ws.Range("H15").Value = "◆開発"
$ws.Hyperlinks.Add($ws.Range("F15"), $ws.Range("F15").Value()) | Out-Null

# Row 16
$ws.Range("A16").Value = "2025-12-08 18:28:28"
$ws.Range("B16").Value = "【フルリモート】WordPressサイトの構築・運用"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Range("F16").Value = "https://www.lancers.jp/work/detail/5449760"
$ws.Range("G16").Value = 58
$ws.Range("H16").Value = "◇サイト ○WordPress"
$ws.Hyperlinks.Add($ws.Range("F16"), $ws.Range("F16").Value()) | Out-Null

# Row 17
$ws.Range("A17").Value = "2025-12-08 18:28:28"
$ws.Range("B17").Value = "初回 WebアプリのiOSアプリ化+IAPサブスク(2週無料)+申請"
$ws.Range("C17").Value = "システム開発"
$ws.Range("D17").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E17").Value = "期限情報なし"
$ws.Range("F17").Value = "https://www.lancers.jp/work/detail/5449067"
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = "◇アプリ"
$ws.Hyperlinks.Add($ws.Range("F17"), $ws.Range("F17").Value()) | Out-Null

# Row 18
$ws.Range("A18").Value = "2025-12-08 18:28:28"
$ws.Range("B18").Value = "【急募】社内システム保守運用・社内スタッフ教育まで依頼できる方を探しています!"
$ws.Range("C18").Value = "システム開発"
$ws.Range("D18").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E18").Value = "期限情報なし"
$ws.Range("F18").Value = "https://www.lancers.jp/work/detail/5449609"
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = ""
$ws.Hyperlinks.Add($ws.Range("F18"), $ws.Range("F18").Value()) | Out-Null

# Row 19
$ws.Range("A19").Value = "2025-12-08 18:28:28"
$ws.Range("B19").Value = "【急募】Shopifyでの3Dカスタムシミュレーター導入設定依頼"
$ws.Range("C19").Value = "システム開発"
$ws.Range("D19").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E19").Value = "期限情報なし"
$ws.Range("F19").Value = "https://www.lancers.jp/work/detail/5449335"
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = ""
$ws.Hyperlinks.Add($ws.Range("F19"), $ws.Range("F19").Value()) | Out-Null

# Row 20
$ws.Range("A20").Value = "2025-12-08 18:28:28"
$ws.Range("B20").Value = "【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え"
$ws.Range("C20").Value = "システム開発"
$ws.Range("D20").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E20").Value = "期限情報なし"
$ws.Range("F20").Value = "https://www.lancers.jp/work/detail/5443568"
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = ""
$ws.Hyperlinks.Add($ws.Range("F20"), $ws.Range("F20").Value()) | Out-Null

# Row 21
$ws.Range("A21").Value = "2025-12-08 18:28:28"
$ws.Range("B21").Value = "【急募】LINEのLステップ構築をサポートしてくれる方"
$ws.Range("C21").Value = "システム開発"
$ws.Range("D21").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E21").Value = "期限情報なし"
$ws.Range("F21").Value = "https://www.lancers.jp/work/detail/5449657"
$ws.Range("G21").Value = 13
$ws.Range("H21").Value = ""
$ws.Hyperlinks.Add($ws.Range("F21"), $ws.Range("F21").Value()) | Out-Null

# Row 22
$ws.Range("A22").Value = "2025-12-08 18:28:28"
$ws.Range("B22").Value = "初回 【継続案件】AWS上でのLAMP環境構築および保守・運用サポートパートナー募集"
$ws.Range("C22").Value = "システム開発"
$ws.Range("D22").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E22").Value = "期限情報なし"
$ws.Range("F22").Value = "https://www.lancers.jp/work/detail/5449313"
$ws.Range("G22").Value = 13
$ws.Range("H22").Value = ""
$ws.Hyperlinks.Add($ws.Range("F22"), $ws.Range("F22").Value()) | Out-Null

# Clear anything beyond the new data range (defensive, in case of stale rows)
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -gt 22) {
    $ws.Range("A23:H" + $lastRow).Clear()
}
